$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "27.721.98"
$ws.Range("E2").Value2 = "  -0.59%  "
$ws.Range("D3").Value2 = "1.848.95"
$ws.Range("E3").Value2 = "  -1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.013"
$ws.Range("E4").Value2 = "  -2.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "319.60"
$ws.Range("E5").Value2 = "  -1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "1.012"
$ws.Range("E6").Value2 = "  -1.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4317"
$ws.Range("E7").Value2 = "  -2.66%  "
$ws.Range("E8").Value2 = "  -2.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.07354"
$ws.Range("E9").Value2 = "  -1.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.8781"
$ws.Range("E10").Value2 = "  -1.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "21.66"
$ws.Range("E11").Value2 = "  -0.62%  "
$ws.Range("D12").Value2 = "1.858.83"
$ws.Range("E12").Value2 = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "6.738"
$ws.Range("E13").Value2 = "  -0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.448"
$ws.Range("E14").Value2 = "  -2.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.07147"
$ws.Range("E15").Value2 = "  -0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "89.21"
$ws.Range("E16").Value2 = "  +4.87%  "
$ws.Range("E17").Value2 = "  -2.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.000009007"
$ws.Range("E18").Value2 = "  -1.53%  "
$ws.Range("E19").Value2 = "  -1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "15.51"
$ws.Range("E20").Value2 = "  -0.84%  "
$ws.Range("D21").Value2 = "27.728.25"
$ws.Range("E21").Value2 = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "5.221"
$ws.Range("E22").Value2 = "  -2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "11.11"
$ws.Range("E23").Value2 = "  -2.16%  "
$ws.Range("D24").Value2 = "2.077.74"
$ws.Range("E24").Value2 = "  -1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "1.994"
$ws.Range("E25").Value2 = "  -1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "155.40"
$ws.Range("E26").Value2 = "  -2.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "18.68"
$ws.Range("E27").Value2 = "  -1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "2.182"
$ws.Range("E28").Value2 = "  +9.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "5.394"
$ws.Range("E29").Value2 = "  -0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "119.15"
$ws.Range("E30").Value2 = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.08946"
$ws.Range("E31").Value2 = "  -1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.236"
$ws.Range("E32").Value2 = "  -0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.7788"
$ws.Range("E33").Value2 = "  -0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "4.576"
$ws.Range("E34").Value2 = "  -1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "2.914"
$ws.Range("E35").Value2 = "  -3.52%  "
$ws.Range("E36").Value2 = "  -2.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "1.136"
$ws.Range("E37").Value2 = "  -0.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.05358"
$ws.Range("E38").Value2 = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.01979"
$ws.Range("E39").Value2 = "  -0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "7.326"
$ws.Range("E40").Value2 = "  +5.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "2.890"
$ws.Range("E41").Value2 = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.5144"
$ws.Range("E42").Value2 = "  -1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.1690"
$ws.Range("E43").Value2 = "  -0.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "8.835"
$ws.Range("E44").Value2 = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "10.70"
$ws.Range("E45").Value2 = "  -0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "108.84"
$ws.Range("E46").Value2 = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.4794"
$ws.Range("E47").Value2 = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.06475"
$ws.Range("E48").Value2 = "  -2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.696"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.013"
$ws.Range("E50").Value2 = "  -2.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.859"
$ws.Range("E51").Value2 = "  -3.30%  "
